$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")

$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
